$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.914.66'
$ws.Range('E2').Value = '  -0.88%  '
$ws.Range('D3').Value = '3.816.63'
$ws.Range('E3').Value = '  -2.10%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '600.36'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.37%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '169.56'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.46%  '
$ws.Range('D7').Value = '3.816.01'
$ws.Range('E7').Value = '  -2.14%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.531'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.22%  '
$ws.Range('E10').Value = '  -0.73%  '
$ws.Range('E11').Value = '  +1.25%  '
$ws.Range('E12').Value = '  +0.85%  '
$ws.Range('E13').Value = '  +9.65%  '
$ws.Range('E14').Value = '  -0.33%  '
$ws.Range('D15').Value = '4.457.81'
$ws.Range('E15').Value = '  -2.08%  '
$ws.Range('D16').Value = '3.812.20'
$ws.Range('E16').Value = '  -2.69%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '18.66'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.81%  '
$ws.Range('D18').Value = '67.969.72'
$ws.Range('E18').Value = '  -0.72%  '
$ws.Range('E19').Value = '  +0.31%  '
$ws.Range('E20').Value = '  +0.16%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.89'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.03%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '470.50'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.51%  '
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('E24').Value = '  -8.95%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.64'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.15%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.31'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.53%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.24'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.25%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.34'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.32%  '
$ws.Range('E29').Value = '  -0.12%  '
$ws.Range('E30').Value = '  -1.51%  '
$ws.Range('D31').Value = '3.965.52'
$ws.Range('E31').Value = '  -2.08%  '
$ws.Range('E32').Value = '  -1.39%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.28'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.36%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '30.80'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.31%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.37'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.73%  '
$ws.Range('D36').Value = '3.782.49'
$ws.Range('E36').Value = '  -2.34%  '
$ws.Range('E37').Value = '  +1.78%  '
$ws.Range('E38').Value = '  +3.92%  '
$ws.Range('E39').Value = '  +1.06%  '
$ws.Range('E40').Value = '  -1.22%  '
$ws.Range('E41').Value = '  -2.24%  '
$ws.Range('E42').Value = '  -0.11%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.319'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.74%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.81'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.20%  '
$ws.Range('E46').Value = '  -1.56%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '410.98'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.43%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '46.50'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.61%  '
$ws.Range('E49').Value = '  -3.81%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '142.56'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.76%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0359'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.15%  '
